$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix data entry error that assigned trees from plot 41 to plot 56: the
# Date column (column C) for rows 3-13 was incorrectly incremented day by
# day instead of all sharing the correct survey date. Reset them to the
# correct date (45149, i.e. 8/11/2023), matching row 2.
for ($row = 3; $row -le 13; $row++) {
    $ws.Cells.Item($row, 3).Value = 45149
}
